$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.048.77"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.900.63"
$ws.Range("E3").Value = "  +1.48%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'312.54"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.5054"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.09344"
$ws.Range("E9").Value = "  -2.50%  "

$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("D11").Value = "'41.79"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("D13").Value = "'20.76"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").Value = "1.885.89"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").Value = "'7.310"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").Value = "'92.42"

$ws.Range("D19").Value = "'0.06579"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").Value = "'17.82"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'6.213"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").Value = "28.110.60"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").Value = "'2.632"
$ws.Range("E26").Value = "  +3.40%  "

$ws.Range("D27").Value = "2.109.25"
$ws.Range("E27").Value = "  +1.27%  "

$ws.Range("D28").Value = "'20.88"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("D29").Value = "'157.11"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").Value = "'127.12"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "'1.087"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("D32").Value = "'0.1067"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("D33").Value = "'5.618"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").Value = "'3.616"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").Value = "'9.653"
$ws.Range("E35").Value = "  +1.62%  "

$ws.Range("D36").Value = "'0.06638"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").Value = "'0.02420"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("D38").Value = "'0.2174"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("D40").Value = "'1.269"
$ws.Range("E40").Value = "  +7.57%  "

$ws.Range("D41").Value = "'0.6389"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").Value = "'4.993"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").Value = "'11.41"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").Value = "'13.31"
$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("D46").Value = "'0.6000"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("D47").Value = "'3.717"
$ws.Range("E47").Value = "  +1.45%  "

$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'2.029"
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").Value = "'122.44"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").Value = "'1.179"
$ws.Range("E51").Value = "  -1.23%  "
